$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 691, shifting existing rows 691..788 down to 692..789
$ws.Rows.Item(691).Insert()

# Populate the newly inserted row 691 with its values
$ws.Cells.Item(691, 1).Value = 3
$ws.Cells.Item(691, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(691, 3).Value = "Coquimbo"
$ws.Cells.Item(691, 4).Value = 45131
$ws.Cells.Item(691, 5).Value = 5
$ws.Cells.Item(691, 6).Value = 100112037
$ws.Cells.Item(691, 7).Value = "Cebollín"
$ws.Cells.Item(691, 8).Value = "Sin especificar"
$ws.Cells.Item(691, 9).Value = "Primera"
$ws.Cells.Item(691, 10).Value = 130
$ws.Cells.Item(691, 11).Value = 4500
$ws.Cells.Item(691, 12).Value = 4500
$ws.Cells.Item(691, 13).Value = 4500
$ws.Cells.Item(691, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(691, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(691, 16).Value = 125
$ws.Cells.Item(691, 17).Value = 36
$ws.Cells.Item(691, 18).Value = "Hortaliza"
